# Updates the cryptocurrency Price (column D) and Volume(1h) (column E) values
# on the "symbol list" worksheet refresh (GitHub Actions, Thu Feb 16 07:30:24 UTC 2023).
#
# Every cell in columns D/E on this sheet is stored as text (inline/shared string),
# e.g. "321.25" or "7.76%". We prefix each new value with a leading apostrophe so
# Excel keeps writing it as text instead of auto-converting the numeric-/percent-
# looking string into a real number (which would change the cell type and drop
# formatting such as trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.21"
$ws.Range("E2").Value = "'7.64%"
$ws.Range("D3").Value = "'49.11"
$ws.Range("E3").Value = "'17.39%"
$ws.Range("D4").Value = "'5.261"
$ws.Range("E4").Value = "'5.03%"
$ws.Range("D5").Value = "'0.08095"
$ws.Range("E5").Value = "'7.38%"
$ws.Range("D6").Value = "'4.615"
$ws.Range("E6").Value = "'5.52%"
$ws.Range("E7").Value = "'3.28%"
$ws.Range("E8").Value = "'30.96%"
$ws.Range("D9").Value = "'0.1322"
$ws.Range("E9").Value = "'11.68%"
$ws.Range("D10").Value = "'0.1947"
$ws.Range("E10").Value = "'6.46%"
$ws.Range("D11").Value = "'0.09517"
$ws.Range("E11").Value = "'5.91%"
$ws.Range("D12").Value = "'0.04466"
$ws.Range("E12").Value = "'10.81%"
$ws.Range("D13").Value = "'0.1048"
$ws.Range("E13").Value = "'-0.19%"
$ws.Range("D14").Value = "'0.001328"
$ws.Range("E14").Value = "'3.64%"
$ws.Range("D15").Value = "'0.005905"
$ws.Range("E15").Value = "'-0.37%"
$ws.Range("D16").Value = "'3.364"
$ws.Range("E16").Value = "'0.75%"
$ws.Range("D17").Value = "'2.434"
$ws.Range("E17").Value = "'1.37%"
$ws.Range("D18").Value = "'0.3392"
$ws.Range("E18").Value = "'1.89%"
$ws.Range("D19").Value = "'8.227"
$ws.Range("E19").Value = "'-0.82%"
$ws.Range("D20").Value = "'0.1415"
$ws.Range("E20").Value = "'3.21%"
$ws.Range("D22").Value = "'0.04305"
$ws.Range("E23").Value = "'3.33%"
$ws.Range("D24").Value = "'0.004243"
$ws.Range("E24").Value = "'8.53%"
$ws.Range("D25").Value = "'0.0001352"
$ws.Range("E25").Value = "'8.11%"
$ws.Range("D26").Value = "'0.0003546"
$ws.Range("E26").Value = "'-4.78%"
$ws.Range("D38").Value = "'0.02711"
$ws.Range("E38").Value = "'12.43%"
$ws.Range("D39").Value = "'0.05579"
$ws.Range("E39").Value = "'7.07%"
$ws.Range("D40").Value = "'0.006311"
$ws.Range("E40").Value = "'0.13%"
$ws.Range("D41").Value = "'0.007675"
$ws.Range("E41").Value = "'-1.67%"
$ws.Range("D42").Value = "'0.1438"
$ws.Range("E42").Value = "'8.45%"
$ws.Range("D43").Value = "'0.007709"
$ws.Range("E43").Value = "'4.23%"
$ws.Range("E44").Value = "'14.23%"
$ws.Range("D46").Value = "'0.00006994"
$ws.Range("E46").Value = "'6.40%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.13%"
$ws.Range("E48").Value = "'30.14%"
$ws.Range("D49").Value = "'0.004007"
$ws.Range("E49").Value = "'-4.63%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.13%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.13%"
